$d = $word.ActiveDocument

# 1. Title heading (appears twice: top heading and bolded run near end) - replace both
$d.Content.Find.Execute("Play FashionTV Highlife Slot Free " + [char]8211 + " Fashion and Style Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play FashionTV Highlife Free | Slot Game Review", 2)

# 2. "What we like" bullet list items
$d.Content.Find.Execute("Unique fashion-themed design", $true, $false, $false, $false, $false, $true, 1, $false, "Stylish design and theme", 2)
$d.Content.Find.Execute("1,280 different paylines", $true, $false, $false, $false, $false, $true, 1, $false, "Wide range of betting options", 2)
$d.Content.Find.Execute("Special Nudging HP1 feature to boost wins", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting bonus features", 2)
$d.Content.Find.Execute("Compatible with both desktop and mobile devices", $true, $false, $false, $false, $false, $true, 1, $false, "Compatibility with mobile devices", 2)

# 3. "What we don't like" bullet list items
$d.Content.Find.Execute("High volatility may not appeal to some players", $true, $false, $false, $false, $false, $true, 1, $false, "High volatility may not appeal to all players", 2)
$d.Content.Find.Execute("Free Spins feature can be difficult to trigger", $true, $false, $false, $false, $false, $true, 1, $false, "Limited number of paylines", 2)

# 4. Italic summary paragraph near end
$d.Content.Find.Execute("Experience exciting gameplay and special features on FashionTV Highlife Slot. Play free on desktop and mobile devices. Win with bonus features and special nudging function.", $true, $false, $false, $false, $false, $true, 1, $false, "Experience the thrill of FashionTV Highlife and play for free. Read our review for more details.", 2)
